# Update COVID stats workbook ("Update countries & provincias Spain")
# - Refresh the "Datos actualizados" timestamp in A1
# - Refresh the per-country case/death numbers that changed between the
#   11:19 and 12:36 data pulls (country names / row order are unchanged)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 21 de Octubre de 2020 a las 12:36"

# --- Updated country rows ----------------------------------------------
# Row 34 - Rumania
$ws.Range("B34").Value = 191102
$ws.Range("C34").Value = 4848
$ws.Range("D34").Value = 137835
$ws.Range("E34").Value = 47202
$ws.Range("G34").Value = 69
$ws.Range("H34").Value = 6065

# Row 52 - Etiopia
$ws.Range("B52").Value = 91763
$ws.Range("C52").Value = 5596
$ws.Range("D52").Value = 54600
$ws.Range("E52").Value = 35018
$ws.Range("H52").Value = 2145

# Row 53 - Honduras
$ws.Range("B53").Value = 90490
$ws.Range("C53").Value = 0
$ws.Range("D53").Value = 43638
$ws.Range("E53").Value = 45481
$ws.Range("G53").Value = 0
$ws.Range("H53").Value = 1371

# Row 54 - Bielorrusia
$ws.Range("B54").Value = 90232
$ws.Range("C54").Value = 851
$ws.Range("D54").Value = 35930
$ws.Range("E54").Value = 51720
$ws.Range("G54").Value = 6
$ws.Range("H54").Value = 2582

# Row 55 - Venezuela
$ws.Range("B55").Value = 88909
$ws.Range("D55").Value = 80503
$ws.Range("E55").Value = 7469
$ws.Range("H55").Value = 937

# Row 56 - Suiza
$ws.Range("B56").Value = 87644
$ws.Range("D56").Value = 80316
$ws.Range("E56").Value = 6581
$ws.Range("H56").Value = 747

# Row 70 - Libia
$ws.Range("B70").Value = 51625
$ws.Range("C70").Value = 719
$ws.Range("D70").Value = 28440
$ws.Range("E70").Value = 22420
$ws.Range("G70").Value = 19
$ws.Range("H70").Value = 765

# Row 100 - Senegal
$ws.Range("B100").Value = 15484
$ws.Range("C100").Value = 25
$ws.Range("D100").Value = 13975
$ws.Range("E100").Value = 1188
$ws.Range("G100").Value = 1
$ws.Range("H100").Value = 321

# Row 131 - Hong Kong
$ws.Range("B131").Value = 5270
$ws.Range("C131").Value = 8
$ws.Range("D131").Value = 5004

# Row 146 - Guyana
$ws.Range("B146").Value = 3797
$ws.Range("C146").Value = 188
$ws.Range("D146").Value = 1341
$ws.Range("E146").Value = 2409
$ws.Range("H146").Value = 47

# Row 147 - Tailandia
$ws.Range("B147").Value = 3796
$ws.Range("C147").Value = 0
$ws.Range("D147").Value = 2796
$ws.Range("E147").Value = 886
$ws.Range("H147").Value = 114

# Row 148 - Gambia
$ws.Range("B148").Value = 3709
$ws.Range("C148").Value = 9
$ws.Range("D148").Value = 3495
$ws.Range("E148").Value = 155
$ws.Range("H148").Value = 59

# Row 149 - Principado de Andorra
$ws.Range("B149").Value = 3655
$ws.Range("D149").Value = 2658
$ws.Range("E149").Value = 879
$ws.Range("H149").Value = 118

# Row 150 - Letonia
$ws.Range("B150").Value = 3623
$ws.Range("D150").Value = 2273
$ws.Range("E150").Value = 1288
$ws.Range("H150").Value = 62
